# Apply the "school name cleanup" edit to the SPS-AttendanceAreas workbook.
#
# Summary of the change (per the commit's xml diff):
#  - Several school-name / attendance-area labels in the "elementary" sheet
#    (column B) were corrected / re-worded, e.g.
#       "Hay"                -> "John Hay"
#       "Concord Int'l"       -> "Concord Intl"
#       "Gatzert, Leschi, Thurgood Marshall" -> "Bailey Gatzert, Leschi, Thurgood Marshall"
#    (the underlying shared-string table gets re-packed as a side effect,
#    but that's just bookkeeping - the actual content change is the text
#    of these 11 cells).
#  - The sheet's active/selected cell moved from B28 to B20.
#  - The workbook window position metadata (xWindow/yWindow) also changed
#    in the source file; we set it for completeness even though it is a
#    purely cosmetic, non-content attribute.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("elementary")

# --- Column B text corrections -------------------------------------------------
$ws.Range("B2").Value  = "Bailey Gatzert, Leschi, Thurgood Marshall"
$ws.Range("B4").Value  = "Viewlands, Daniel Bagley"
$ws.Range("B6").Value  = "Broadview-Thomson, Northgate"
$ws.Range("B18").Value = "Thurgood Marshall, John Muir, Beacon Hill Intl, Kimball, Leschi"
$ws.Range("B19").Value = "Beacon Hill Intl, Kimball"
$ws.Range("B21").Value = "Dearborn Park, Van Asselt, Martin Luther King Jr., Wing Luke, Dunlap, Hawthorne"
$ws.Range("B25").Value = "Concord Intl"
$ws.Range("B26").Value = "Green Lake, Laurelhurst"
$ws.Range("B27").Value = "John Hay, Coe"
$ws.Range("B28").Value = "John Hay"
$ws.Range("B29").Value = "B. F. Day, Green Lake"

# --- Move the active selection from B28 to B20 ---------------------------------
$ws.Activate()
$ws.Range("B20").Select()

# --- Cosmetic workbook window position (best effort) ----------------------------
$win = $wb.Windows.Item(1)
$win.Left = 7700
$win.Top = 900
